$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.NumberFormat = "General"
    $cell.ClearFormats()
}

Set-TextValue $ws.Cells.Item(2, 3) "199"
Set-TextValue $ws.Cells.Item(2, 4) "474016.00"
Set-TextValue $ws.Cells.Item(3, 3) "1064"
Set-TextValue $ws.Cells.Item(3, 4) "3500837.33"
Set-TextValue $ws.Cells.Item(4, 3) "432"
Set-TextValue $ws.Cells.Item(4, 4) "1827442.36"
Set-TextValue $ws.Cells.Item(5, 3) "121"
Set-TextValue $ws.Cells.Item(5, 4) "580128.09"
Set-TextValue $ws.Cells.Item(8, 3) "46"
Set-TextValue $ws.Cells.Item(8, 4) "95000.00"
Set-TextValue $ws.Cells.Item(9, 3) "63"
Set-TextValue $ws.Cells.Item(9, 4) "164607.56"
Set-TextValue $ws.Cells.Item(10, 3) "376"
Set-TextValue $ws.Cells.Item(10, 4) "1377261.71"
Set-TextValue $ws.Cells.Item(11, 3) "150"
Set-TextValue $ws.Cells.Item(11, 4) "653601.77"
Set-TextValue $ws.Cells.Item(15, 3) "106"
Set-TextValue $ws.Cells.Item(15, 4) "300126.38"
Set-TextValue $ws.Cells.Item(16, 3) "440"
Set-TextValue $ws.Cells.Item(16, 4) "1409153.23"
Set-TextValue $ws.Cells.Item(17, 3) "162"
Set-TextValue $ws.Cells.Item(17, 4) "701553.96"
Set-TextValue $ws.Cells.Item(19, 3) "18"
Set-TextValue $ws.Cells.Item(19, 4) "121216.00"
Set-TextValue $ws.Cells.Item(20, 3) "19"
Set-TextValue $ws.Cells.Item(20, 4) "40621.00"
Set-TextValue $ws.Cells.Item(21, 3) "54"
Set-TextValue $ws.Cells.Item(21, 4) "130800.00"
Set-TextValue $ws.Cells.Item(33, 3) "115"
Set-TextValue $ws.Cells.Item(33, 4) "339173.00"
Set-TextValue $ws.Cells.Item(34, 3) "606"
Set-TextValue $ws.Cells.Item(34, 4) "2048578.10"
Set-TextValue $ws.Cells.Item(35, 3) "257"
Set-TextValue $ws.Cells.Item(35, 4) "1369531.11"
Set-TextValue $ws.Cells.Item(36, 3) "82"
Set-TextValue $ws.Cells.Item(36, 4) "476394.00"
Set-TextValue $ws.Cells.Item(37, 3) "27"
Set-TextValue $ws.Cells.Item(37, 4) "180500.00"
Set-TextValue $ws.Cells.Item(38, 3) "26"
Set-TextValue $ws.Cells.Item(38, 4) "56200.00"
Set-TextValue $ws.Cells.Item(39, 3) "47"
Set-TextValue $ws.Cells.Item(39, 4) "133886.00"
Set-TextValue $ws.Cells.Item(40, 3) "186"
Set-TextValue $ws.Cells.Item(40, 4) "521991.00"
Set-TextValue $ws.Cells.Item(41, 3) "88"
Set-TextValue $ws.Cells.Item(41, 4) "320429.00"
Set-TextValue $ws.Cells.Item(43, 3) "13"
Set-TextValue $ws.Cells.Item(43, 4) "65500.00"
Set-TextValue $ws.Cells.Item(44, 3) "74"
Set-TextValue $ws.Cells.Item(44, 4) "185183.00"
Set-TextValue $ws.Cells.Item(45, 3) "33"
Set-TextValue $ws.Cells.Item(45, 4) "114357.84"
Set-TextValue $ws.Cells.Item(46, 3) "102"
Set-TextValue $ws.Cells.Item(46, 4) "459974.61"
Set-TextValue $ws.Cells.Item(50, 3) "18"
Set-TextValue $ws.Cells.Item(50, 4) "39850.00"
Set-TextValue $ws.Cells.Item(51, 3) "111"
Set-TextValue $ws.Cells.Item(51, 4) "335698.17"
Set-TextValue $ws.Cells.Item(52, 3) "644"
Set-TextValue $ws.Cells.Item(52, 4) "2485548.42"
Set-TextValue $ws.Cells.Item(53, 3) "281"
Set-TextValue $ws.Cells.Item(53, 4) "1330578.76"
Set-TextValue $ws.Cells.Item(54, 3) "102"
Set-TextValue $ws.Cells.Item(54, 4) "653274.23"
Set-TextValue $ws.Cells.Item(55, 3) "29"
Set-TextValue $ws.Cells.Item(55, 4) "173213.00"
Set-TextValue $ws.Cells.Item(57, 3) "798"
Set-TextValue $ws.Cells.Item(57, 4) "2127800.24"
Set-TextValue $ws.Cells.Item(58, 3) "3928"
Set-TextValue $ws.Cells.Item(58, 4) "12997226.28"
Set-TextValue $ws.Cells.Item(59, 3) "2033"
Set-TextValue $ws.Cells.Item(59, 4) "9190618.89"
Set-TextValue $ws.Cells.Item(60, 3) "705"
Set-TextValue $ws.Cells.Item(60, 4) "3707842.28"
Set-TextValue $ws.Cells.Item(61, 3) "154"
Set-TextValue $ws.Cells.Item(61, 4) "1124223.00"
Set-TextValue $ws.Cells.Item(62, 3) "4"
Set-TextValue $ws.Cells.Item(62, 4) "40000.00"
Set-TextValue $ws.Cells.Item(63, 3) "388"
Set-TextValue $ws.Cells.Item(63, 4) "973422.27"
Set-TextValue $ws.Cells.Item(64, 3) "36"
Set-TextValue $ws.Cells.Item(64, 4) "97961.00"
Set-TextValue $ws.Cells.Item(65, 3) "141"
Set-TextValue $ws.Cells.Item(65, 4) "360837.69"
Set-TextValue $ws.Cells.Item(66, 3) "68"
Set-TextValue $ws.Cells.Item(66, 4) "229421.21"
Set-TextValue $ws.Cells.Item(67, 3) "23"
Set-TextValue $ws.Cells.Item(67, 4) "102876.00"
Set-TextValue $ws.Cells.Item(68, 3) "5"
Set-TextValue $ws.Cells.Item(68, 4) "22500.00"
Set-TextValue $ws.Cells.Item(69, 3) "9"
Set-TextValue $ws.Cells.Item(69, 4) "19500.00"
Set-TextValue $ws.Cells.Item(83, 3) "236"
Set-TextValue $ws.Cells.Item(83, 4) "605326.09"
Set-TextValue $ws.Cells.Item(84, 3) "924"
Set-TextValue $ws.Cells.Item(84, 4) "3012146.26"
Set-TextValue $ws.Cells.Item(85, 3) "343"
Set-TextValue $ws.Cells.Item(85, 4) "1410985.70"
Set-TextValue $ws.Cells.Item(86, 3) "124"
Set-TextValue $ws.Cells.Item(86, 4) "628484.52"
Set-TextValue $ws.Cells.Item(96, 3) "429"
Set-TextValue $ws.Cells.Item(96, 4) "1343302.36"
Set-TextValue $ws.Cells.Item(97, 3) "182"
Set-TextValue $ws.Cells.Item(97, 4) "740304.72"
Set-TextValue $ws.Cells.Item(98, 3) "62"
Set-TextValue $ws.Cells.Item(98, 4) "331911.73"
